# "Draw fittings with blocks"
# Update O/P (and dependent) columns on all three LEONARDO sheets: apply the
# "0" integer number format to the Q,resa/Q,tot cells, refresh a couple of
# power formulas, resize the header row on the 3.0 PLUS sheet, and leave the
# selections where the author left them.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# LEONARDO 5.5
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LEONARDO 5.5")

$ws1.Range("O3").NumberFormat = "0"
$ws1.Range("P3").NumberFormat = "0"

$ws1.Range("O4").NumberFormat = "0"
$ws1.Range("O4").Formula = "=51.8*D4*E4"

$ws1.Range("P4").NumberFormat = "0"
$ws1.Range("P4").Formula = "=51.8*E4*1.1*D4"

# ---------------------------------------------------------------------------
# LEONARDO 3.5
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LEONARDO 3.5")

$ws2.Range("O3").NumberFormat = "0"
$ws2.Range("P3").NumberFormat = "0"
$ws2.Range("O4").NumberFormat = "0"
$ws2.Range("P4").NumberFormat = "0"

# ---------------------------------------------------------------------------
# LEONARDO 3.0 PLUS
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("LEONARDO 3.0 PLUS")

$ws3.Rows.Item(2).RowHeight = 42

$ws3.Range("O3").NumberFormat = "0"

$ws3.Range("P3").NumberFormat = "0"
$ws3.Range("P3").Formula = "=82.3*E3*1.1*D3"

$ws3.Range("O4").NumberFormat = "0"
$ws3.Range("O4").Formula = "=80.9*E4*D4"

$ws3.Range("P4").NumberFormat = "0"

# ---------------------------------------------------------------------------
# Restore per-sheet selections / scroll position (order matters: the last
# sheet activated ends up as the workbook's active tab, which must stay on
# "LEONARDO 3.0 PLUS" to match the saved file).
# ---------------------------------------------------------------------------
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 2
$ws1.Range("O3:P4").Select()

$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 4
$ws2.Range("P3").Select()

$ws3.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 6
$ws3.Range("P4").Select()
